# Min Price Rows Page Test
# Adds a new "testMinRowSelection" test row to the RunManager sheet and a
# matching data row to the TestData sheet, and flips the now-superseded
# checkBrokenLinks "Execution Flag" from yes -> no on both sheets.

$wb = $excel.ActiveWorkbook

$runMgr = $wb.Worksheets.Item("RunManager")
$testData = $wb.Worksheets.Item("TestData")

# ---------------------------------------------------------------------
# RunManager sheet: flip the checkBrokenLinks flag off, then append the
# new testMinRowSelection row (A8:E8).
# ---------------------------------------------------------------------
$runMgr.Range("C7").Value = "no"

$runMgr.Range("A8").Value = "testMinRowSelection"
$runMgr.Range("B8").Value = "To Test the Min Row is Selected"
$runMgr.Range("C8").Value = "yes"

# D8/E8 hold the text "1" (same as D2:E7) rather than a number, so copy an
# existing text "1" cell instead of assigning a literal that Excel would
# auto-convert to a numeric value.
$runMgr.Range("D2").Copy() | Out-Null
$runMgr.Range("D8").PasteSpecial() | Out-Null
$runMgr.Range("E2").Copy() | Out-Null
$runMgr.Range("E8").PasteSpecial() | Out-Null

$runMgr.Range("A8").Select() | Out-Null

# ---------------------------------------------------------------------
# TestData sheet: flip checkBrokenLinks's flag off, then append the new
# testMinRowSelection data row (A12:G12).
# ---------------------------------------------------------------------
$testData.Range("B11").Value = "no"

$testData.Range("A12").Value = "testMinRowSelection"
$testData.Range("B12").Value = "yes"
$testData.Range("C12").Value = "chrome"

# D12:G12 all hold a lone apostrophe placeholder, same as D11:G11. A bare
# "'" assigned through .Value is interpreted by Excel as a text-prefix
# marker (producing an empty cell), so copy the existing placeholder cells
# instead.
$testData.Range("D11").Copy() | Out-Null
$testData.Range("D12").PasteSpecial() | Out-Null
$testData.Range("E11").Copy() | Out-Null
$testData.Range("E12").PasteSpecial() | Out-Null
$testData.Range("F11").Copy() | Out-Null
$testData.Range("F12").PasteSpecial() | Out-Null
$testData.Range("G11").Copy() | Out-Null
$testData.Range("G12").PasteSpecial() | Out-Null

$testData.Range("A11").Select() | Out-Null
